$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false,
                             $true, 1, $false, $replace, 2) | Out-Null
}

# 1. Heading3: "Industrial Processes I" -> "Industrial Chemical Processes"
Replace-Text "Industrial Processes I" "Industrial Chemical Processes"

# 2. Ativação date
Replace-Text "Ativação: 01/01/2024" "Ativação: 01/01/2025"

# 3. Objetivos (PT)
Replace-Text "Conferir aos alunos uma visão geral da indústria química e correlatas, processos e produtos, e uma visão global das matérias primas mais importantes da indústria química." "Conferir aos alunos uma visão geral da indústria química e correlatas, bem como das principais características dos processos desta indústria."

# 4. Objectives (EN)
Replace-Text "Objectives:Check the students an overview of the chemical industry and related industries, processes and products, and an overview of the most important raw materials in the chemical industry.." "Providing to the students an overview of the chemical and related industries, as well as the main features of the processes and production arrangements of this industry."

# 5. Programa (PT) - long text (>255 chars), Find.Execute search text is capped at 255
#    chars by Word, so set the whole (single-run) paragraph's Range.Text directly.
$p14 = $d.Paragraphs.Item(14)
$p14.Range.Text = "O conteúdo desta disciplina será de acordo com os tópicos a serem programados, devendo abordar assuntos relevantes relacionados a processos químicos e correlatas."

# 6. Programa (EN) - also long text
$p15 = $d.Paragraphs.Item(15)
$p15.Range.Text = "The content of this subject will be in accordance with the topics to be programmed, and must address relevant subjects related to chemical and related processes."

# 7. Método
Replace-Text "Aulas expositivas, filmes e leituras de artigos técnicos" "Aulas expositivas, desenvolvimento de exercícios em sala e fora de sala de aula, discussão de casos práticos."

# 8. Critério
Replace-Text "Provas em sala, apresentações em sala, entrega de exercícios ou casos práticos elaborados fora de sala de aula e frequência." "A nota (NOTA) será composta por uma destas opções: prova em sala, apresentações em sala, entrega de exercícios ou casos práticos elaborados fora de sala de aula. A estas opções será incorporado, para cada aluno, seu respectivo percentual de frequência no cálculo da nota final (NF), conforme a fórmula explicitada abaixo:NF = NOTA x % FREQ."

# 9. Norma de recuperação
Replace-Text "Prova escrita para alunos que tenham média final maior ou igual a 3,0 (Três) e inferior a 5,0 (Cinco). A nota final será a média aritmética entre a média final e a prova escrita." "Frequência mínima de 70% e nota igual ou superior a 3,00 e inferior a 5,00 possibilita aplicação de prova escrita de recuperação valendo 10,00 pontos."

# 10. Bibliografia - single run paragraph, long text on both sides, set directly.
#     Note: the kept prefix (through "...c1997.Revistas:") contains U+037E (a
#     Greek question mark that is visually a semicolon) at two spots in the
#     source document, after "Weinheim " and after "R. Norris" - reproduce it
#     via [char] so the byte-for-byte original text is preserved exactly.
$gq = [char]0x037E
$p19 = $d.Paragraphs.Item(19)
$p19.Range.Text = "Ullmann’s encyclopedia of industrial chemistry; Editorial advisory board, Giuseppe Bellussi et al.; 7th, completely revised edition; Weinheim $($gq) New York : WileyVCH, 2011.Encyclopedia of Chemical Processing; Edited by Sunggyu Lee; New York : Taylor & Francis, 2006.Kirk, Raymond Eller. Encyclopedia of chemical technology / Herman F.Mark et al. New York: John Wiley, 1984.Manual Econômico da Indústria Química - MEIQ / Centro de Pesquisas e Desenvolvimento; 8ed; Camaçari: CEPED, 2007.Shreve, R. Norris$($gq) BRINK JR., J. A. Indústrias de processos químicos. Tradução de Horácio Macedo; 4.ed. Rio de Janeiro: Editora Guanabara Koogan, 2008, c1997.Revistas:Brazilian Journal of Chemical Engineering, São Paulo, SP: Brazilian Society of Chemical Engineering, v. 11, n. 1, 1995-;"
